{"js": "// Update the date title and every equation cell in the practice table.\n// The table's cell COUNT never changes (20 rows x 5 cols before AND after);\n// the diff's apparent \"insert a cell / delete a cell\" in the last data row\n// is really just every cell value from that point on shifting by one\n// position, so we simply rewrite all 25 equation cells (in row-major,\n// left-to-right order) plus the title paragraph with their new values.\n\nconst titleNew = \"2023-12-23 Saturday\";\n\n// New equation values, in the same left-to-right / top-to-bottom order as\n// the 5 data rows (row 0, 4, 8, 12, 16) of the table.\nconst newValues = [\n  [\"63\u00f72=\", \"87\u00f78=\", \"43\u00f75=\", \"60\u00f76=\", \"62\u00f78=\"],\n  [\"60\u00f75=\", \"96\u00f75=\", \"68\u00f72=\", \"81\u00f76=\", \"61\u00f78=\"],\n  [\"90\u00f76=\", \"67\u00f77=\", \"32\u00f79=\", \"62\u00f78=\", \"75\u00f78=\"],\n  [\"14\u00f75=\", \"43\u00f78=\", \"53\u00f76=\", \"52\u00f72=\", \"78\u00f77=\"],\n  [\"47\u00f78=\", \"89\u00f73=\", \"35\u00f73=\", \"81\u00f76=\", \"76\u00f74=\"],\n];\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2023-12-22 Friday\" -> \"2023-12-23 Saturday\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(titleNew, \"Replace\");\n\n// 2) Table cells: walk the rows that actually contain text (every 4th row,\n// starting at row 0) and overwrite each cell's value in place.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet dataRowIdx = 0;\nfor (let r = 0; r < rows.items.length; r++) {\n  const cells = rows.items[r].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  if (cells.items.length === 0 || cells.items[0].value === \"\") continue;\n\n  const vals = newValues[dataRowIdx];\n  for (let c = 0; c < cells.items.length && c < vals.length; c++) {\n    cells.items[c].value = vals[c];\n  }\n  dataRowIdx++;\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and every equation cell in the practice table.\n# The table's cell COUNT never changes (20 rows x 5 cols before AND after);\n# the diff's apparent \"insert a cell / delete a cell\" in the last data row\n# is really just every cell value from that point on shifting by one\n# position, so we simply rewrite all 25 equation cells (in row-major,\n# left-to-right order) plus the title paragraph with their new values.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2023-12-22 Friday\" -> \"2023-12-23 Saturday\".\n$d.Paragraphs.Item(1).Range.Text = \"2023-12-23 Saturday\"\n\n# 2) Table cells: the table has 20 rows x 5 columns, but only every 4th row\n# (1, 5, 9, 13, 17 in 1-based COM indexing) actually holds an equation; the\n# rows between are blank spacer rows. Overwrite each of those 5 data rows.\n$newValues = @(\n    @(\"63\u00f72=\", \"87\u00f78=\", \"43\u00f75=\", \"60\u00f76=\", \"62\u00f78=\"),\n    @(\"60\u00f75=\", \"96\u00f75=\", \"68\u00f72=\", \"81\u00f76=\", \"61\u00f78=\"),\n    @(\"90\u00f76=\", \"67\u00f77=\", \"32\u00f79=\", \"62\u00f78=\", \"75\u00f78=\"),\n    @(\"14\u00f75=\", \"43\u00f78=\", \"53\u00f76=\", \"52\u00f72=\", \"78\u00f77=\"),\n    @(\"47\u00f78=\", \"89\u00f73=\", \"35\u00f73=\", \"81\u00f76=\", \"76\u00f74=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($i = 0; $i -lt $newValues.Count; $i++) {\n    $rowIndex = 1 + 4 * $i\n    $row = $newValues[$i]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $t.Cell($rowIndex, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
